$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2275641025641026
$ws.Range("C2").Value = 0.5160256410256411
$ws.Range("J2").Value = 0.009615384615384616
$ws.Range("O2").Value = 0.003205128205128205
$ws.Range("P2").Value = 0.1666666666666667
$ws.Range("S2").Value = 0.07692307692307693
$ws.Range("C3").Value = 0.05780346820809248
$ws.Range("J3").Value = 0.0115606936416185
$ws.Range("P3").Value = 0.7514450867052023
$ws.Range("S3").Value = 0.1791907514450867
$ws.Range("J4").Value = 0.06382978723404255
$ws.Range("P4").Value = 0.851063829787234
$ws.Range("S4").Value = 0.0851063829787234
$ws.Range("B6").Value = 0.05531914893617021
$ws.Range("D6").Value = 0.008510638297872341
$ws.Range("F6").Value = 0.08085106382978724
$ws.Range("J6").Value = 0.3021276595744681
$ws.Range("O6").Value = 0.01702127659574468
$ws.Range("Q6").Value = 0.148936170212766
$ws.Range("R6").Value = 0.05106382978723404
$ws.Range("S6").Value = 0.3361702127659574
$ws.Range("B7").Value = 0.08839779005524862
$ws.Range("D7").Value = 0.01657458563535912
$ws.Range("F7").Value = 0.06077348066298342
$ws.Range("J7").Value = 0.1546961325966851
$ws.Range("O7").Value = 0.01657458563535912
$ws.Range("Q7").Value = 0.1878453038674033
$ws.Range("R7").Value = 0.07734806629834254
$ws.Range("S7").Value = 0.3977900552486188
$ws.Range("B8").Value = 0.0931174089068826
$ws.Range("D8").Value = 0.02226720647773279
$ws.Range("F8").Value = 0.05668016194331984
$ws.Range("J8").Value = 0.145748987854251
$ws.Range("O8").Value = 0.01821862348178137
$ws.Range("Q8").Value = 0.1761133603238866
$ws.Range("R8").Value = 0.08502024291497975
$ws.Range("S8").Value = 0.402834008097166
$ws.Range("B9").Value = 0.09142857142857143
$ws.Range("D9").Value = 0.01142857142857143
$ws.Range("F9").Value = 0.02857142857142857
$ws.Range("J9").Value = 0.1371428571428571
$ws.Range("O9").Value = 0.01142857142857143
$ws.Range("Q9").Value = 0.2342857142857143
$ws.Range("R9").Value = 0.09142857142857143
$ws.Range("S9").Value = 0.3942857142857143
$ws.Range("B10").Value = 0.1075268817204301
$ws.Range("D10").Value = 0.02150537634408602
$ws.Range("E10").Value = 0.0007168458781362007
$ws.Range("F10").Value = 0.07670250896057347
$ws.Range("J10").Value = 0.1541218637992832
$ws.Range("O10").Value = 0.008602150537634409
$ws.Range("Q10").Value = 0.2200716845878136
$ws.Range("R10").Value = 0.08172043010752689
$ws.Range("S10").Value = 0.3290322580645161
$ws.Range("G11").Value = 0.1263157894736842
$ws.Range("J11").Value = 0.112280701754386
$ws.Range("K11").Value = 0.1929824561403509
$ws.Range("L11").Value = 0.5614035087719298
$ws.Range("S11").Value = 0.007017543859649123
$ws.Range("G12").Value = 0.7393939393939394
$ws.Range("J12").Value = 0.2060606060606061
$ws.Range("K12").Value = 0.02424242424242424
$ws.Range("L12").Value = 0.01818181818181818
$ws.Range("S12").Value = 0.01212121212121212
$ws.Range("F15").Value = 0.01785714285714286
$ws.Range("I15").Value = 0.08482142857142858
$ws.Range("J15").Value = 0.3794642857142857
$ws.Range("K15").Value = 0.1071428571428571
$ws.Range("M15").Value = 0.01339285714285714
$ws.Range("O15").Value = 0.0625
$ws.Range("S15").Value = 0.2098214285714286
$ws.Range("F16").Value = 0.01401869158878505
$ws.Range("H16").Value = 0.2242990654205607
$ws.Range("I16").Value = 0.07476635514018691
$ws.Range("J16").Value = 0.3785046728971962
$ws.Range("K16").Value = 0.102803738317757
$ws.Range("M16").Value = 0.01869158878504673
$ws.Range("O16").Value = 0.0514018691588785
$ws.Range("S16").Value = 0.1355140186915888
$ws.Range("F17").Value = 0.01789264413518887
$ws.Range("H17").Value = 0.1848906560636183
$ws.Range("I17").Value = 0.07554671968190854
$ws.Range("J17").Value = 0.4274353876739563
$ws.Range("K17").Value = 0.06759443339960239
$ws.Range("M17").Value = 0.02385685884691849
$ws.Range("O17").Value = 0.07157057654075547
$ws.Range("S17").Value = 0.1312127236580517
$ws.Range("F18").Value = 0.02010050251256281
$ws.Range("H18").Value = 0.1658291457286432
$ws.Range("I18").Value = 0.05527638190954774
$ws.Range("J18").Value = 0.4271356783919598
$ws.Range("K18").Value = 0.1005025125628141
$ws.Range("M18").Value = 0.01507537688442211
$ws.Range("O18").Value = 0.06030150753768844
$ws.Range("S18").Value = 0.1557788944723618
$ws.Range("F19").Value = 0.01658767772511848
$ws.Range("H19").Value = 0.2353870458135861
$ws.Range("I19").Value = 0.0703001579778831
$ws.Range("J19").Value = 0.3609794628751975
$ws.Range("K19").Value = 0.1003159557661927
$ws.Range("M19").Value = 0.01737756714060032
$ws.Range("N19").Value = 0.0007898894154818325
$ws.Range("O19").Value = 0.07109004739336493
$ws.Range("S19").Value = 0.127172195892575
